$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVT")

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Drillthrough"
$ws.Range("C10").Value = "1.Generate a chart with some data`n2.Create new report page and in DrillThrough add the fields for drillthrough.`n3. Right click on the chart, select the Drillthrough option from the menu. "
$ws.Range("D10").Value = "1. On right click of the chart and selecting the drillthrough option from the context-menu , the report will drillthrough to the newly created report page."

$ws.Range("C10:E10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 75

$ws.Range("C11:C29").Select() | Out-Null
